$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: mark numeric-looking price cells as Text so Excel keeps them as literal strings
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Step 2: write the new cell values
$ws.Range('D2').Value = '30.452.72'
$ws.Range('E2').Value = '  +0.98%  '
$ws.Range('D3').Value = '1.878.18'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '247.11'
$ws.Range('E5').Value = '  +5.48%  '
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('D7').Value = '0.4766'
$ws.Range('E7').Value = '  +1.83%  '
$ws.Range('D8').Value = '0.2898'
$ws.Range('E8').Value = '  +1.23%  '
$ws.Range('D9').Value = '0.06526'
$ws.Range('E9').Value = '  +0.56%  '
$ws.Range('D10').Value = '21.83'
$ws.Range('E10').Value = '  +3.31%  '
$ws.Range('E11').Value = '  -0.18%  '
$ws.Range('D12').Value = '96.80'
$ws.Range('E12').Value = '  +3.18%  '
$ws.Range('D13').Value = '0.7370'
$ws.Range('E13').Value = '  +7.79%  '
$ws.Range('D14').Value = '1.879.17'
$ws.Range('E14').Value = '  +0.79%  '
$ws.Range('D15').Value = '5.126'
$ws.Range('E15').Value = '  +1.42%  '
$ws.Range('D16').Value = '273.17'
$ws.Range('E16').Value = '  +1.63%  '
$ws.Range('D17').Value = '30.477.97'
$ws.Range('E17').Value = '  +1.08%  '
$ws.Range('E18').Value = '  +2.34%  '
$ws.Range('D19').Value = '0.000007605'
$ws.Range('E19').Value = '  -0.45%  '
$ws.Range('E20').Value = '  -0.04%  '
$ws.Range('D21').Value = '2.130.63'
$ws.Range('E21').Value = '  +0.37%  '
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('D23').Value = '5.240'
$ws.Range('E23').Value = '  +1.70%  '
$ws.Range('D24').Value = '6.181'
$ws.Range('E24').Value = '  +1.28%  '
$ws.Range('D25').Value = '9.305'
$ws.Range('E25').Value = '  -0.44%  '
$ws.Range('D26').Value = '164.22'
$ws.Range('E26').Value = '  -0.88%  '
$ws.Range('D27').Value = '18.89'
$ws.Range('E27').Value = '  +1.84%  '
$ws.Range('D28').Value = '1.941'
$ws.Range('E28').Value = '  +2.59%  '
$ws.Range('D29').Value = '1.373'
$ws.Range('E29').Value = '  +0.58%  '
$ws.Range('D30').Value = '0.09965'
$ws.Range('E30').Value = '  +0.41%  '
$ws.Range('E31').Value = '  +4.52%  '
$ws.Range('D32').Value = '4.317'
$ws.Range('E32').Value = '  +2.19%  '
$ws.Range('E33').Value = '  +1.53%  '
$ws.Range('D34').Value = '0.04775'
$ws.Range('E34').Value = '  +1.86%  '
$ws.Range('D35').Value = '1.124'
$ws.Range('E35').Value = '  +0.54%  '
$ws.Range('D36').Value = '0.7004'
$ws.Range('E36').Value = '  +1.57%  '
$ws.Range('D37').Value = '2.716'
$ws.Range('E37').Value = '  +0.30%  '
$ws.Range('E38').Value = '  +1.81%  '
$ws.Range('D39').Value = '2.732'
$ws.Range('E39').Value = '  -0.93%  '
$ws.Range('D40').Value = '6.331'
$ws.Range('E40').Value = '  -0.39%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = '70.59'
$ws.Range('E41').Value = '  -0.95%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value = '1.939'
$ws.Range('E42').Value = '  +2.12%  '
$ws.Range('D43').Value = '0.4187'
$ws.Range('E43').Value = '  +3.10%  '
$ws.Range('D44').Value = '0.9999'
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('D45').Value = '0.8381'
$ws.Range('E45').Value = '  +0.48%  '
$ws.Range('D46').Value = '102.67'
$ws.Range('E46').Value = '  +0.43%  '
$ws.Range('D47').Value = '9.272'
$ws.Range('E47').Value = '  +1.58%  '
$ws.Range('D48').Value = '7.087'
$ws.Range('E48').Value = '  +1.86%  '
$ws.Range('D49').Value = '35.52'
$ws.Range('E49').Value = '  +4.16%  '
$ws.Range('D50').Value = '925.98'
$ws.Range('E50').Value = '  -0.83%  '
$ws.Range('E51').Value = '  +1.10%  '

# Step 3: restore default (Normal) style on the cells we forced to Text,
# so the saved file does not carry a lingering style index on them
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
